# Saldo.xlsx edit script
# Applies the account-balance refresh described in the commit:
#   - the top account (ELAINE/29000) is replaced by five much larger
#     balances (FABRICIO, ANA, SUELY, NATAL, CEZAR)
#   - the old "ANA 477.7" row is removed (its account now appears higher
#     up with its corrected balance)
#   - a "RAPHAELA" row is inserted further down with an updated balance,
#     and the stale RAPHAELA (218.85) / SUELY (200) rows that duplicated
#     those accounts further down the list are removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Work from the bottom of the sheet upward so every row number below refers
# to the ORIGINAL (unedited) layout - earlier (higher-numbered) edits never
# disturb the row numbers used by later (lower-numbered) ones.
# ---------------------------------------------------------------------------

# 1) Remove the stale "SUELY / 004216401 / 200" row (row 87) - this account's
#    balance was corrected and moved near the top of the sheet.
$ws.Rows.Item(87).Delete()

# 2) Remove the stale "RAPHAELA / 005366255 / 218.85" row (row 80) - replaced
#    by a fresh row with the corrected balance inserted below.
$ws.Rows.Item(80).Delete()

# 3) Insert the corrected RAPHAELA row just above row 71 (SERGIO).
$ws.Rows.Item(71).Insert()
$ws.Cells.Item(71,1).NumberFormat = "@"
$ws.Cells.Item(71,1).Value = "005366255"
$ws.Cells.Item(71,2).Value = "RAPHAELA"
$ws.Cells.Item(71,3).Value = 245.44

# 4) Remove the stale "ANA / 004261201 / 477.7" row (row 33) - replaced by the
#    corrected balance inserted near the top of the sheet.
$ws.Rows.Item(33).Delete()

# 5) Remove the original top data row ("ELAINE / 005018038 / 29000", row 2).
$ws.Rows.Item(2).Delete()

# 6) Insert the five corrected high-value rows at the top of the data block.
$ws.Rows.Item(2).Resize(5).Insert()

$topRows = @(
    @("004570632", "FABRICIO", 51069.05),
    @("004261201", "ANA",      22790.39),
    @("004216401", "SUELY",    17768.16),
    @("001759765", "NATAL",    14419.09),
    @("004482090", "CEZAR",     5034.16)
)

for ($i = 0; $i -lt $topRows.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r,1).NumberFormat = "@"
    $ws.Cells.Item($r,1).Value = $topRows[$i][0]
    $ws.Cells.Item($r,2).Value = $topRows[$i][1]
    $ws.Cells.Item($r,3).Value = $topRows[$i][2]
}
